# Insert a new (blank) column before column N ("Late") on the "Repayment
# schedule" sheet. This shifts the old N ("Late") -> O, the old spacer O -> P,
# and the old P ("Outstanding") -> Q, matching the new layout that adds a
# blank spacer column ahead of "Late" (mirroring the existing spacer ahead of
# "Outstanding").
$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Repayment schedule")
$ws3.Columns.Item(14).Insert()

# Make "Repayment schedule" the active sheet/tab, with the given cell
# selected, matching the saved view state in the workbook.
$ws3.Activate()
$ws3.Range("L18").Select() | Out-Null
